# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Rows 19-22 also got re-ranked (Avalanche / WrappedliquidstakedEther2.0 /
# Uniswap / ShibaInu swapped positions), so those rows' Coin/Link columns
# are rewritten too, not just Price/Volume(1h).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use Text format for D:E while writing so numeric-looking strings
# (e.g. "1.001", "0.4667") are stored as literal text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.281.44"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.886.93"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "238.15"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.4667"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("D8").Value = "0.2823"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").Value = "0.06571"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Value = "19.76"
$ws.Range("E10").Value = "  +5.07%  "
$ws.Range("D11").Value = "0.07771"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "97.93"
$ws.Range("E12").Value = "  -3.29%  "
$ws.Range("D13").Value = "1.886.12"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "5.120"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "0.6670"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "281.92"
$ws.Range("E16").Value = "  +9.72%  "
$ws.Range("D17").Value = "30.280.42"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "12.61"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.133.29"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.359"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.000007298"
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "6.167"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "9.344"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "165.61"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").Value = "19.11"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "1.989"
$ws.Range("E28").Value = "  -3.20%  "
$ws.Range("D29").Value = "1.377"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "0.09756"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").Value = "4.455"
$ws.Range("E31").Value = "  -5.18%  "
$ws.Range("D32").Value = "1.483"
$ws.Range("D33").Value = "4.166"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").Value = "0.04694"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "0.7068"
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("D36").Value = "1.095"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "0.01868"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").Value = "6.657"
$ws.Range("E39").Value = "  +6.79%  "
$ws.Range("D40").Value = "2.524"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").Value = "72.09"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("D42").Value = "0.8685"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "104.07"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").Value = "0.4200"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "989.75"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "7.206"
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "9.263"
$ws.Range("E49").Value = "  +4.82%  "
$ws.Range("E50").Value = "  -2.98%  "
$ws.Range("E51").Value = "  -2.39%  "

# Restore default styling on the touched columns (values already committed as text).
$ws.Range("D2:E51").Style = "Normal"

